# Add a "Currency" column to the "Stock log" sheet (after "Operation", before
# "Number of Shares"), populate existing rows with "USD", and add a new
# "ENBRIDGE INC" (ENB.TO) holding priced in CAD. Mirror the new holding into
# the "Portfolio Summary" sheet and refresh the "Total Return" totals.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: "Stock log" - insert the new "Currency" column (D) by re-writing
# columns D..M one slot to the right (E..N), then fill in the new D column.
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Stock log")

# Header row: rewrite columns D..M (old) into E..N (new), then set D1.
# N1 falls beyond the sheet's previous used range, so it starts out with no
# style - copy the bold/centered/bordered header formatting over from C1
# (any existing header cell) before writing its text.
$ws1.Cells.Item(1, 3).Copy()
$ws1.Cells.Item(1, 14).PasteSpecial(-4122)

$headers = @("Number of Shares", "Buy Price", "Market Price", "Annual Dividend per Share", "Date of Purchase", "Cost Basis", "Market Value", "Capital Gains", "Dividends Paid", "Total Return")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws1.Cells.Item(1, 5 + $i).Value = $headers[$i]
}
$ws1.Cells.Item(1, 4).Value = "Currency"

# Data rows 2-5: shift existing D..M values one column right into E..N, and
# set the new Currency column (D) to "USD" for all of them.
$rows = @(
    @(687, 23.55, 21.69, 1.78, "22.06.2023", 16178.85, 14901.03, -7.9, 611.4300000000001, -4.28),
    @(400, 24.66, 21.69, 1.78, "27.07.2023", 9864, 8676, -12.04, 356, -8.75),
    @(983, 75.65000000000001, 76.3, 2.66, "27.07.2023", 74363.95000000001, 75002.89999999999, 0.86, 1376.2, 2.76),
    @(-100, 23, 21.69, 1.78, "25.11.2023", -2300, -2169, -5.7, -44, -3.86)
)

for ($r = 0; $r -lt $rows.Length; $r++) {
    $rowNum = $r + 2
    $data = $rows[$r]
    for ($i = 0; $i -lt $data.Length; $i++) {
        $ws1.Cells.Item($rowNum, 5 + $i).Value = $data[$i]
    }
    $ws1.Cells.Item($rowNum, 4).Value = "USD"
}

# New row 6: ENBRIDGE INC / ENB.TO, bought in CAD.
$ws1.Cells.Item(6, 1).Value = "ENBRIDGE INC"
$ws1.Cells.Item(6, 2).Value = "ENB.TO"
$ws1.Cells.Item(6, 3).Value = "Buy"
$ws1.Cells.Item(6, 4).Value = "CAD"
$ws1.Cells.Item(6, 5).Value = 100
$ws1.Cells.Item(6, 6).Value = 45
$ws1.Cells.Item(6, 7).Value = 35.7825
$ws1.Cells.Item(6, 8).Value = 3.66
$ws1.Cells.Item(6, 9).Value = "20.05.2023"
$ws1.Cells.Item(6, 10).Value = 3375
$ws1.Cells.Item(6, 11).Value = 3578.25
$ws1.Cells.Item(6, 12).Value = 6.02
$ws1.Cells.Item(6, 13).Value = 133.5
$ws1.Cells.Item(6, 14).Value = 10.39

# ---------------------------------------------------------------------
# Sheet 2: "Portfolio Summary" - append the ENB.TO summary row.
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Portfolio Summary")

$ws2.Cells.Item(4, 1).Value = "ENB.TO"
# Copy the bold/centered/bordered header-style formatting used by the other
# ticker cells in column A (e.g. A3) onto the new A4 cell.
$ws2.Cells.Item(3, 1).Copy()
$ws2.Cells.Item(4, 1).PasteSpecial(-4122)

$ws2.Cells.Item(4, 2).Value = 100
$ws2.Cells.Item(4, 3).Value = 3375
$ws2.Cells.Item(4, 4).Value = 3578.25
$ws2.Cells.Item(4, 5).Value = 133.5
$ws2.Cells.Item(4, 6).Value = 33.75
$ws2.Cells.Item(4, 7).Value = 6.01
$ws2.Cells.Item(4, 8).Value = 10.36

# ---------------------------------------------------------------------
# Sheet 3: "Total Return" - refresh "My Portfolio" totals now that the new
# CAD holding has been converted/included.
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Total Return")

$ws3.Cells.Item(2, 3).Value = -1.47
$ws3.Cells.Item(2, 4).Value = 0.95
